# Add JAS Dec 2017 entries (journal_URL.xlsx)
# Appends 18 new rows (154-171) to Sheet1, reusing the existing lookup
# values (journal/class/type) and adding the new article URLs, with the
# first 5 new URL cells (154-158) turned into real hyperlinks like the
# existing styled rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (journal, year, month, class, type, URL)
$newRows = @(
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5290"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5420"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5430"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5439"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5447"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5455"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5466"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5474"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5485"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5507"),
    @("journal of animal science", 2017, 12, "monogastric", "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5516"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5547"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5563"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5573"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5584"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5597"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5606"),
    @("journal of animal science", 2017, 12, "ruminant",    "original article", "https://www.animalsciencepublications.org/publications/jas/articles/95/12/5617")
)

$startRow = 154
$hyperlinkCount = 5

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
    $ws.Range("E$r").Value = $vals[4]
    $ws.Range("F$r").Value = $vals[5]

    if ($i -lt $hyperlinkCount) {
        $ws.Hyperlinks.Add($ws.Range("F$r"), $vals[5])
    }
}

# Re-apply the sheet's standard hyperlink cell style (as used by every
# other hyperlinked URL cell, e.g. F2) to the newly linked cells so they
# reuse the workbook's existing hyperlink style instead of a fresh one.
$lastHyperlinkRow = $startRow + $hyperlinkCount - 1
$ws.Range("F2").Copy()
$ws.Range("F$($startRow):F$($lastHyperlinkRow)").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final cursor position/selection on the sheet.
$ws.Range("D170").Select()

Write-Output "Added rows $startRow to $($startRow + $newRows.Count - 1); hyperlinks on F$startRow:F$lastHyperlinkRow"
